$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns K1:Q1 (additional turn parameters)
$ws.Range("K1").Value = "Average G Force"
$ws.Range("L1").Value = "Maximum G Force"
$ws.Range("M1").Value = "Average Slope"
$ws.Range("N1").Value = "Maximum Slope"
$ws.Range("O1").Value = "G Force at Maximum Slope"
$ws.Range("P1").Value = "Turn Radius"
$ws.Range("Q1").Value = "Turn Type"

# Sample data for the new Turn Type column
$ws.Range("Q2").Value = "Carving"
$ws.Range("Q3").Value = "Skidding"

# Update the view to match the committed selection/scroll position
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("Q4").Select()
